# Add results from latest run: extend the nowcast table with a new
# "2025-08-30" row (row 12) and refresh all previously-computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged text, just rewritten to keep shared-string ---
# --- bookkeeping consistent with a fresh data pull) ---------------------
$headers = @("Row", "Prognose", "surveys", "production", "orders", "turnover", "financial", "labor market", "prices", "national accounts", "Revision")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Row labels (dates), column A, rows 2-12 -----------------------------
# These are stored as literal text (not real dates) in the source data, so
# force a Text format before assigning, then restore the default ("Normal")
# style so the cell doesn't pick up a stray custom number format.
$dates = @("2025-03-30", "2025-04-15", "2025-04-30", "2025-05-15", "2025-05-30", "2025-06-15", "2025-06-30", "2025-07-15", "2025-07-30", "2025-08-15", "2025-08-30")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $cell = $ws.Cells.Item($i + 2, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.Style = "Normal"
}

# --- Numeric data, columns B-K, rows 2-12 --------------------------------
# (kept as strings and cast to [double] so exact decimal/exponent literals
# round-trip losslessly through the interpreter)
$data = @(
    @("0.29709854136322567", "0", "0", "0", "0", "0", "0", "0", "0", "0"),
    @("0.33062244563092913", "0", "0.0030438001602629498", "-0.0014021649252704584", "-0.00020963998092435925", "0.0095543603729728206", "0.00036541805271730953", "0.023164970937489839", "0", "-0.00099284034954461031"),
    @("0.33117442231183869", "-0.005125543865749066", "0", "-2.9146367000719994e-05", "-0.000784837835961609", "0", "0.00023098125196061218", "0.01885046545394075", "-0.012840360626762913", "0.00025041867048247735"),
    @("0.28966893304103769", "-0.023477831965147258", "-0.010728351226332082", "-0.0021349455197368556", "0.00051170828875677157", "-0.021087488531141771", "0.0005681892815415973", "0.012454465131502203", "0", "0.0023887652697563988"),
    @("0.46856888134199209", "0.17755516705908869", "0", "0.0065574170035253712", "-0.00096532528568271506", "0", "0.00052875320353447082", "-0.018402685053069644", "0", "0.01362662137355819"),
    @("0.46872320259272643", "0", "-0.0078615843573403774", "0.0021598879453898381", "0.00057449854748437652", "0.00129920767633342", "0", "-0.0015137631739424043", "0", "0.0054960746128094873"),
    @("0.48816075435993178", "0.055670379043488032", "0", "-0.010654414368616172", "-0.018214597016723569", "0", "0.00011848456575567823", "-0.012801524768589182", "0", "0.0053192243118905802"),
    @("0.11031598546064891", "0", "-0.28155991367988059", "0.016213041319337599", "-0.067825131685867235", "-0.042338439166350314", "-0.0029286853241497816", "0.0009912560505876823", "0", "-0.00039689641296020151"),
    @("0.13533364330204709", "0.016319618193590665", "0", "0.0053743168987946989", "-0.0063772385943842089", "0", "-0.00085379975121073861", "0.028827708610119178", "-0.001866641060319452", "-0.016406306455191977"),
    @("0.35863962355102252", "0", "0.38832920448094133", "-0.034566561302180748", "-0.030448521635134396", "-0.057458722341774085", "-0.0047435493307956995", "0.043484208099153911", "0", "-0.081290077721234894"),
    @("0.35419746423945359", "0.03279689678936154", "0", "-0.0062880362889653077", "-0.002179940370633317", "0", "-0.00078220666027488471", "-0.02853893168935347", "0", "0.00055005890829651349")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = [double]$row[$c]
    }
}

# --- Column E width reverts to the narrower default used by the other ---
# --- numeric columns once the new column was added. ----------------------
$ws.Columns.Item(5).ColumnWidth = 15.24609375
